# Update the GitHub repository link shown in the "REPOSITORIO GIT" table
# cell from the old SSH-style remote URL to the new HTTPS project URL.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "git@github.com:cmercadoloayza/DigitalHarbor-Clinica.git",  # FindText
    $true,                                                      # MatchCase
    $false,                                                     # MatchWholeWord
    $false,                                                     # MatchWildcards
    $false,                                                     # MatchSoundsLike
    $false,                                                     # MatchAllWordForms
    $true,                                                      # Forward
    1,                                                           # Wrap (wdFindContinue)
    $false,                                                     # Format
    "https://github.com/cmercadoloayza/Proyecto-Clinica",       # ReplaceWith
    2                                                            # Replace (wdReplaceAll)
)
